$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# FE sheet: remove all existing rows (header + 3 data rows) -> sheet left empty
# ---------------------------------------------------------------
$wsFE = $wb.Worksheets.Item("FE")
$wsFE.Rows("1:4").Delete() | Out-Null
$wsFE.Rows("1:1").Select() | Out-Null

# ---------------------------------------------------------------
# SE sheet: add six new exam rows (5-10) below the existing blank styled row 4
# ---------------------------------------------------------------
$wsSE = $wb.Worksheets.Item("SE")

$wsSE.Range("A5").Value = "'DM"
$wsSE.Range("B5").Value = "'Friday 09/07/2021"
$wsSE.Range("C5").Value = "'09:00"
$wsSE.Range("D5").Value = "'11:00"

$wsSE.Range("A6").Value = "'MPMC"
$wsSE.Range("B6").Value = "'Friday 09/07/2021"
$wsSE.Range("C6").Value = "'02:00"
$wsSE.Range("D6").Value = "'04:00"

$wsSE.Range("A7").Value = "'FL & AT"
$wsSE.Range("B7").Value = "'Saturday 10/07/2021"
$wsSE.Range("C7").Value = "'09:00"
$wsSE.Range("D7").Value = "'11:00"

$wsSE.Range("A8").Value = "'MADF"
$wsSE.Range("B8").Value = "'Saturday 10/07/2021"
$wsSE.Range("C8").Value = "'02:00"
$wsSE.Range("D8").Value = "'04:00"

$wsSE.Range("A9").Value = "'OOSE"
$wsSE.Range("B9").Value = "'Monday 12/07/2021"
$wsSE.Range("C9").Value = "'09:00"
$wsSE.Range("D9").Value = "'11:00"

$wsSE.Range("A10").Value = "'ECO"
$wsSE.Range("B10").Value = "'Monday 12/07/2021"
$wsSE.Range("C10").Value = "'02:00"
$wsSE.Range("D10").Value = "'04:00"

$wsSE.Rows("5:10").Select() | Out-Null

# ---------------------------------------------------------------
# BE sheet: selection cursor moved from D2 to the whole row 2
# ---------------------------------------------------------------
$wsBE = $wb.Worksheets.Item("BE")
$wsBE.Rows("2:2").Select() | Out-Null

# ---------------------------------------------------------------
# DT sheet: row 6 updated (FE/IT1/2 -> SE/Internal Test 3/6), row 7 removed
# ---------------------------------------------------------------
$wsDT = $wb.Worksheets.Item("DT")

$wsDT.Range("A6").Value = "'SE"
$wsDT.Range("C6").Value = "'Internal Test 3"
$wsDT.Range("D6").Value = "'6"

$wsDT.Rows("7:7").Delete() | Out-Null
$wsDT.Rows("6:6").Select() | Out-Null
$wsDT.Activate() | Out-Null
